# Apply the edits described by the commit "null/del filter works, but json wont be created"
# The change replaces several per-article free-text fields (manufacturer type, short
# description, manufacturer name, supplier article number) with status/marker texts
# describing validation problems found for each sample row, updates the BRUTTO PREIS
# (gross price) test values, and fills in two previously empty cells on the last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: article removed ("Artikel entfernt")
$ws.Range("C2").Value = "Artikel entfernt"
$ws.Range("D2").Value = "Artikel entfernt"
$ws.Range("G2").Value = "Artikel entfernt"
$ws.Range("P2").Value = "Artikel entfernt"
$ws.Range("S2").Value = 100

# Row 3: AX number missing ("AX Nummer fehlt")
$ws.Range("C3").Value = "AX Nummer fehlt"
$ws.Range("D3").Value = "AX Nummer fehlt"
$ws.Range("G3").Value = "AX Nummer fehlt"
$ws.Range("P3").Value = "AX Nummer fehlt"
$ws.Range("S3").Value = 200

# Row 4: price unit missing ("PE fehlt")
$ws.Range("C4").Value = "PE fehlt"
$ws.Range("D4").Value = "PE fehlt"
$ws.Range("G4").Value = "PE fehlt"
$ws.Range("P4").Value = "PE fehlt"
$ws.Range("S4").Value = 300

# Row 5: price is empty ("Preis ist leer")
$ws.Range("C5").Value = "Preis ist leer"
$ws.Range("D5").Value = "Preis ist leer"
$ws.Range("G5").Value = "Preis ist leer"
$ws.Range("P5").Value = "Preis ist leer"
$ws.Range("S5").Value = 400

# Row 6: everything fits ("ALLES PASST"), fill in previously empty EAN and LAND cells.
# E6 takes the same EAN text as E5 ("7612270054053"); copy/paste-special values from
# E5 so the digit string is stored as text (matching the EAN column elsewhere) instead
# of being auto-converted to a number.
$ws.Range("C6").Value = "ALLES PASST"
$ws.Range("D6").Value = "ALLES PASST"
$ws.Range("E5").Copy()
$ws.Range("E6").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("G6").Value = "ALLES PASST"
$ws.Range("P6").Value = "ALLES PASST"
$ws.Range("Q6").Value = "CH"
$ws.Range("S6").Value = 500

# Update the active selection shown when the sheet was last saved
$ws.Range("P5").Select()
